$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Convert the "OFFERS" column (D2:D14) from text values to real numbers,
#    fix the one value that was wrong (Aliens Group: 12 -> 112), and apply
#    the bordered "Normal" style that Excel uses for plain numeric cells.
$values = @(1, 7, 1, 8, 2, 2, 2, 14, 23, 112, 9, 6, 16)

for ($i = 0; $i -lt $values.Length; $i++) {
  $r = $i + 2
  $cell = $ws.Cells.Item($r, 4)
  $cell.Value = $values[$i]
  $cell.Style = "Normal"
  $cell.Borders.LineStyle = 1
  $cell.Borders.Color = 0
  $cell.Borders.Weight = -4138
}

# 2. Shrink the AutoFilter / filter-database range so it no longer spans
#    the OFFERS column (now B1:C14 instead of B1:D14).
$ws.AutoFilterMode = $false
$ws.Range("B1:C14").AutoFilter()

foreach ($n in $wb.Names) {
  if ($n.Name -eq "Sheet1!_FilterDatabase") {
    $n.RefersTo = "=Sheet1!`$B`$1:`$C`$14"
  }
}

Write-Host "done"
